$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# --- Fix a couple of pre-existing data-quality typos in the stock sheet ---
# Stray full-width space inside company names.
$ws.Cells.Item(2, 2).Value = "中美矽晶製品股份有限公司"
$ws.Cells.Item(3, 2).Value = "聯成化學科技股份有限公司"

# Malformed number (mixed half-width/full-width commas) -> plain digits.
$ws.Cells.Item(6, 7).Value = "1508000"

# --- New column: property_category (issue #5: stock data output to json file) ---
# Insert a blank column before the existing "date" column (column H) so the
# trailing columns (date / legislator_name / legislator_id) shift right by
# one, inheriting their original formatting automatically.
$ws.Columns.Item(8).Insert()

$ws.Cells.Item(1, 8).Value = "property_category"
$ws.Cells.Item(2, 8).Value = "stock"
$ws.Cells.Item(3, 8).Value = "stock"
$ws.Cells.Item(4, 8).Value = "stock"
$ws.Cells.Item(5, 8).Value = "stock"
$ws.Cells.Item(6, 8).Value = "stock"
